$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Início")
$ws2 = $wb.Worksheets.Item("Cad_Empresa")

# Insert a new row at row 7 on "Início" (this shifts the "GESTOR RESPONSÁVEL"
# label and everything below it down by one row, matching the target layout).
$ws1.Rows.Item(7).Insert()

# --- Cad_Empresa: fill in sample data -------------------------------------
# Company name (consumed by formulas on "Início", so it must be set before
# those formulas are written).
$ws2.Range("B6").Value = "transportadora rápida web"

# --- Início: add the three instructional notes near the bottom -----------
# (order chosen so the shared-string table matches the target layout)
$ws1.Range("B22").Value = "Pegar o valor tudo maiúsculo: digita o sinal =maiuscula / digita tab/ e clica na céçuça que tem o valor desejado/dá enter"
$ws1.Range("B21").Value = "Pegar o mesmo conteúdo de texto da célula: clica na célula/insere o sinal = e clica na outra célula/enter"
$ws1.Range("B23").Value = "Copiar as primeiras letras maiusculas: insere o sinal = juntamente com pri.maiuscula/ tab/clica no conteúdo/enter"

# --- Cad_Empresa: remaining sample data -----------------------------------
$ws2.Range("B10").Value = "maria brown silva"
$ws2.Range("B18").Value = 1234

# --- Início: formulas referencing Cad_Empresa -----------------------------
$ws1.Range("B6").Formula = "=Cad_Empresa!B6"

# B7 reuses B6's formatting (style) before getting its own formula.
$ws1.Range("B6").Copy()
$ws1.Range("B7").PasteSpecial(-4122)
$ws1.Range("B7").Formula = "=UPPER(Cad_Empresa!B6)"

$ws1.Range("B10").Formula = "=PROPER(Cad_Empresa!B10)"

# --- View state: active cell / selection on each sheet --------------------
$ws2.Activate() | Out-Null
$ws2.Range("B14").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("C14").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
